$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel doesn't
# auto-convert them to numbers (losing trailing zeros / exponent form).
$forceTextCells = @("D5", "D6", "D7", "D9", "D11", "D12", "D14", "D16", "D19", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}


# Row 2
$ws.Range("D2").Value = '79.745.70'
$ws.Range("E2").Value = '  +4.71%  '

# Row 3
$ws.Range("D3").Value = '3.216.65'
$ws.Range("E3").Value = '  +6.09%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '210.66'
$ws.Range("E5").Value = '  +6.97%  '

# Row 6
$ws.Range("D6").Value = '640.82'
$ws.Range("E6").Value = '  +3.41%  '

# Row 7
$ws.Range("D7").Value = '0.261'
$ws.Range("E7").Value = '  +27.78%  '

# Row 8
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").Value = '0.602'
$ws.Range("E9").Value = '  +9.86%  '

# Row 10
$ws.Range("D10").Value = '3.213.48'
$ws.Range("E10").Value = '  +6.08%  '

# Row 11
$ws.Range("D11").Value = '0.603'
$ws.Range("E11").Value = '  +37.03%  '

# Row 12
$ws.Range("D12").Value = '0.0000264'
$ws.Range("E12").Value = '  +38.29%  '

# Row 13
$ws.Range("E13").Value = '  +3.52%  '

# Row 14
$ws.Range("D14").Value = '5.43'
$ws.Range("E14").Value = '  +3.67%  '

# Row 15
$ws.Range("D15").Value = '3.807.04'
$ws.Range("E15").Value = '  +6.11%  '

# Row 16
$ws.Range("D16").Value = '32.52'
$ws.Range("E16").Value = '  +12.46%  '

# Row 17
$ws.Range("D17").Value = '79.566.64'
$ws.Range("E17").Value = '  +4.55%  '

# Row 18
$ws.Range("D18").Value = '3.210.95'
$ws.Range("E18").Value = '  +6.02%  '

# Row 19
$ws.Range("D19").Value = '14.68'
$ws.Range("E19").Value = '  +9.03%  '

# Row 20
$ws.Range("E20").Value = '  +28.46%  '

# Row 21
$ws.Range("D21").Value = '9.40'

# Row 22
$ws.Range("D22").Value = '446.89'
$ws.Range("E22").Value = '  +17.21%  '

# Row 23
$ws.Range("D23").Value = '5.27'
$ws.Range("E23").Value = '  +20.34%  '

# Row 24
$ws.Range("D24").Value = '4.85'
$ws.Range("E24").Value = '  +12.57%  '

# Row 26
$ws.Range("D26").Value = '77.58'
$ws.Range("E26").Value = '  +7.22%  '

# Row 27
$ws.Range("D27").Value = '10.95'
$ws.Range("E27").Value = '  +12.29%  '

# Row 28
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.44%  '

# Row 29
$ws.Range("D29").Value = '0.0000126'
$ws.Range("E29").Value = '  +17.20%  '

# Row 30
$ws.Range("D30").Value = '9.24'
$ws.Range("E30").Value = '  +11.83%  '

# Row 31
$ws.Range("E31").Value = '  +0.54%  '

# Row 32
$ws.Range("D32").Value = '563.71'
$ws.Range("E32").Value = '  +14.25%  '

# Row 33
$ws.Range("D33").Value = '1.53'
$ws.Range("E33").Value = '  +10.12%  '

# Row 34
$ws.Range("E34").Value = '  +32.32%  '

# Row 35
$ws.Range("E35").Value = '  +6.63%  '

# Row 36
$ws.Range("D36").Value = '23.23'
$ws.Range("E36").Value = '  +13.00%  '

# Row 37
$ws.Range("D37").Value = '0.122'
$ws.Range("E37").Value = '  +18.56%  '

# Row 38
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.09%  '

# Row 39
$ws.Range("D39").Value = '0.415'
$ws.Range("E39").Value = '  +9.77%  '

# Row 40
$ws.Range("D40").Value = '163.07'
$ws.Range("E40").Value = '  +0.54%  '

# Row 41
$ws.Range("D41").Value = '20.28'
$ws.Range("E41").Value = '  +1.30%  '

# Row 42
$ws.Range("D42").Value = '5.74'
$ws.Range("E42").Value = '  +12.63%  '

# Row 43
$ws.Range("D43").Value = '194.88'
$ws.Range("E43").Value = '  +2.20%  '

# Row 45
$ws.Range("D45").Value = '1.84'
$ws.Range("E45").Value = '  +12.52%  '

# Row 46
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '2.71'
$ws.Range("E46").Value = '  +12.36%  '

# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.803'
$ws.Range("E47").Value = '  +3.67%  '

# Row 48
$ws.Range("D48").Value = '1.35'
$ws.Range("E48").Value = '  +8.58%  '

# Row 49
$ws.Range("D49").Value = '43.00'
$ws.Range("E49").Value = '  +4.01%  '

# Row 50
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = '4.32'
$ws.Range("E50").Value = '  +11.56%  '

# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '26.00'
$ws.Range("E51").Value = '  +17.62%  '

# Reset style on the force-text cells back to Normal (no explicit style)
# now that the text values are safely stored, to match original formatting.
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
